$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Cases")
$ws.Cells.Item(23, 11).Value = 101
$ws.Cells.Item(23, 28).Value = 3213
$ws.Cells.Item(24, 11).Value = 128
$ws.Cells.Item(24, 28).Value = 4299
$ws.Cells.Item(25, 11).Value = 182
$ws.Cells.Item(25, 28).Value = 5476
$ws.Cells.Item(26, 11).Value = 226
$ws.Cells.Item(26, 28).Value = 6592
$ws.Cells.Item(27, 11).Value = 258
$ws.Cells.Item(27, 28).Value = 7327
$ws.Cells.Item(28, 11).Value = 284
$ws.Cells.Item(28, 28).Value = 8029
$ws.Cells.Item(29, 11).Value = 300
$ws.Cells.Item(29, 28).Value = 9341
$ws.Cells.Item(30, 11).Value = 343
$ws.Cells.Item(30, 28).Value = 10319
$ws.Cells.Item(31, 11).Value = 393
$ws.Cells.Item(31, 28).Value = 11495
$ws.Cells.Item(32, 11).Value = 433
$ws.Cells.Item(32, 28).Value = 12718
$ws.Cells.Item(33, 11).Value = 474
$ws.Cells.Item(33, 28).Value = 14006
$ws.Cells.Item(34, 11).Value = 495
$ws.Cells.Item(34, 28).Value = 14888
$ws.Cells.Item(35, 11).Value = 519
$ws.Cells.Item(35, 28).Value = 15652
$ws.Cells.Item(36, 11).Value = 535
$ws.Cells.Item(36, 28).Value = 16651
$ws.Cells.Item(37, 11).Value = 547
$ws.Cells.Item(37, 28).Value = 17634
$ws.Cells.Item(38, 11).Value = 592
$ws.Cells.Item(38, 28).Value = 18672
$ws.Cells.Item(39, 11).Value = 622
$ws.Cells.Item(39, 28).Value = 19734
$ws.Cells.Item(40, 11).Value = 649
$ws.Cells.Item(40, 28).Value = 20652
$ws.Cells.Item(41, 11).Value = 657
$ws.Cells.Item(41, 28).Value = 21266
$ws.Cells.Item(42, 11).Value = 668
$ws.Cells.Item(42, 28).Value = 21686
$ws.Cells.Item(43, 11).Value = 679
$ws.Cells.Item(43, 28).Value = 22346
$ws.Cells.Item(44, 11).Value = 684
$ws.Cells.Item(44, 28).Value = 22996
$ws.Cells.Item(45, 11).Value = 698
$ws.Cells.Item(45, 28).Value = 23649
$ws.Cells.Item(46, 11).Value = 719
$ws.Cells.Item(46, 28).Value = 24322
$ws.Cells.Item(47, 11).Value = 730
$ws.Cells.Item(47, 28).Value = 24815
$ws.Cells.Item(48, 11).Value = 743
$ws.Cells.Item(48, 28).Value = 25281
$ws.Cells.Item(49, 11).Value = 747
$ws.Cells.Item(49, 28).Value = 25599
$ws.Cells.Item(50, 11).Value = 750
$ws.Cells.Item(50, 28).Value = 25839
$ws.Cells.Item(51, 11).Value = 753
$ws.Cells.Item(51, 28).Value = 26147
$ws.Cells.Item(52, 11).Value = 759
$ws.Cells.Item(52, 28).Value = 26467
$ws.Cells.Item(53, 11).Value = 764
$ws.Cells.Item(53, 28).Value = 26780
$ws.Cells.Item(54, 11).Value = 770
$ws.Cells.Item(54, 28).Value = 27088
$ws.Cells.Item(55, 11).Value = 773
$ws.Cells.Item(55, 28).Value = 27413
$ws.Cells.Item(56, 11).Value = 778
$ws.Cells.Item(56, 28).Value = 27599
$ws.Cells.Item(57, 11).Value = 785
$ws.Cells.Item(57, 28).Value = 27789
$ws.Cells.Item(58, 11).Value = 787
$ws.Cells.Item(58, 28).Value = 27948
$ws.Cells.Item(59, 11).Value = 791
$ws.Cells.Item(59, 28).Value = 28152

$ws = $wb.Worksheets.Item("Fatalities")
$ws.Cells.Item(23, 11).Value = 1
$ws.Cells.Item(23, 28).Value = 37
$ws.Cells.Item(25, 11).Value = 2
$ws.Cells.Item(25, 28).Value = 57
$ws.Cells.Item(27, 11).Value = 4
$ws.Cells.Item(27, 28).Value = 98
$ws.Cells.Item(28, 11).Value = 4
$ws.Cells.Item(28, 28).Value = 113
$ws.Cells.Item(29, 11).Value = 5
$ws.Cells.Item(29, 28).Value = 143
$ws.Cells.Item(31, 11).Value = 8
$ws.Cells.Item(31, 28).Value = 200
$ws.Cells.Item(32, 11).Value = 8
$ws.Cells.Item(32, 28).Value = 247
$ws.Cells.Item(34, 11).Value = 12
$ws.Cells.Item(34, 28).Value = 335
$ws.Cells.Item(35, 11).Value = 15
$ws.Cells.Item(35, 28).Value = 378
$ws.Cells.Item(36, 11).Value = 19
$ws.Cells.Item(36, 28).Value = 435
$ws.Cells.Item(37, 11).Value = 23
$ws.Cells.Item(37, 28).Value = 502
$ws.Cells.Item(38, 11).Value = 26
$ws.Cells.Item(38, 28).Value = 564
$ws.Cells.Item(39, 11).Value = 28
$ws.Cells.Item(39, 28).Value = 625
$ws.Cells.Item(40, 11).Value = 30
$ws.Cells.Item(40, 28).Value = 686
$ws.Cells.Item(41, 11).Value = 30
$ws.Cells.Item(41, 28).Value = 749
$ws.Cells.Item(42, 11).Value = 33
$ws.Cells.Item(42, 28).Value = 802
$ws.Cells.Item(43, 11).Value = 34
$ws.Cells.Item(43, 28).Value = 858
$ws.Cells.Item(44, 11).Value = 35
$ws.Cells.Item(44, 28).Value = 917
$ws.Cells.Item(46, 11).Value = 36
$ws.Cells.Item(46, 28).Value = 1034
$ws.Cells.Item(47, 11).Value = 36
$ws.Cells.Item(47, 28).Value = 1084
$ws.Cells.Item(48, 11).Value = 36
$ws.Cells.Item(48, 28).Value = 1120
$ws.Cells.Item(49, 11).Value = 36
$ws.Cells.Item(49, 28).Value = 1174
$ws.Cells.Item(50, 11).Value = 37
$ws.Cells.Item(50, 28).Value = 1204
$ws.Cells.Item(54, 11).Value = 40
$ws.Cells.Item(54, 28).Value = 1380
$ws.Cells.Item(55, 11).Value = 40
$ws.Cells.Item(55, 28).Value = 1419
$ws.Cells.Item(56, 11).Value = 41
$ws.Cells.Item(56, 28).Value = 1440
$ws.Cells.Item(57, 11).Value = 42
$ws.Cells.Item(57, 28).Value = 1479
$ws.Cells.Item(58, 11).Value = 43
$ws.Cells.Item(58, 28).Value = 1514

$ws = $wb.Worksheets.Item("Hospitalized")
$ws.Cells.Item(23, 11).Value = 9
$ws.Cells.Item(23, 28).Value = 422
$ws.Cells.Item(29, 11).Value = 29
$ws.Cells.Item(29, 28).Value = 1206
$ws.Cells.Item(30, 11).Value = 43
$ws.Cells.Item(30, 28).Value = 1364
$ws.Cells.Item(31, 11).Value = 45
$ws.Cells.Item(31, 28).Value = 1462
$ws.Cells.Item(32, 11).Value = 52
$ws.Cells.Item(32, 28).Value = 1612
$ws.Cells.Item(33, 11).Value = 58
$ws.Cells.Item(33, 28).Value = 1792
$ws.Cells.Item(34, 11).Value = 63
$ws.Cells.Item(34, 28).Value = 1879
$ws.Cells.Item(35, 11).Value = 63
$ws.Cells.Item(35, 28).Value = 1996
$ws.Cells.Item(36, 11).Value = 58
$ws.Cells.Item(36, 28).Value = 2174
$ws.Cells.Item(38, 11).Value = 59
$ws.Cells.Item(38, 28).Value = 2284
$ws.Cells.Item(39, 11).Value = 60
$ws.Cells.Item(39, 28).Value = 2348
$ws.Cells.Item(40, 11).Value = 51
$ws.Cells.Item(40, 28).Value = 2325
$ws.Cells.Item(41, 11).Value = 52
$ws.Cells.Item(41, 28).Value = 2314
$ws.Cells.Item(54, 11).Value = 33
$ws.Cells.Item(54, 28).Value = 1582
$ws.Cells.Item(55, 11).Value = 29
$ws.Cells.Item(55, 28).Value = 1538
